# The "Förändrad" (Changed) date column (C) for rows 2-52 was bulk-updated
# from the Excel date serial 45192 (2023-09-23) to 45202 (2023-10-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C52").Value = 45202
